$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 235.15384
$ws.Range("I33").Value = 111
$ws.Range("J33").Value = 433.8
$ws.Range("K33").Value = 111
$ws.Range("L33").Value = 433.8
$ws.Range("M33").Value = 118
$ws.Range("N33").Value = -891.8
$ws.Range("H97").Value = 900
$ws.Range("J97").Value = 900
$ws.Range("L97").Value = 2700
$ws.Range("N97").Value = -3692
$ws.Range("H98").Value = 843.8
$ws.Range("I98").Value = 830.6667
$ws.Range("J98").Value = 896.3333
$ws.Range("K98").Value = 830.6667
$ws.Range("L98").Value = 896.3333
$ws.Range("M98").Value = 667.3333
$ws.Range("N98").Value = -3892.3333
$ws.Range("H122").Value = 843.8
$ws.Range("I122").Value = 830.6667
$ws.Range("J122").Value = 896.3333
$ws.Range("K122").Value = 2492.0001
$ws.Range("L122").Value = 2688.9999
$ws.Range("M122").Value = -42.0001000000002
$ws.Range("N122").Value = -7588.9999
$ws.Range("H135").Value = 31504.06
$ws.Range("I135").Value = 36754.215
$ws.Range("J135").Value = 2103.2
$ws.Range("K135").Value = 330787.9349999999
$ws.Range("L135").Value = 18928.8
$ws.Range("M135").Value = -328252.9349999999
$ws.Range("N135").Value = -23998.8
$ws.Range("H137").Value = 3848197.5
$ws.Range("J137").Value = 4250
$ws.Range("L137").Value = 12750
$ws.Range("N137").Value = -17850
$ws.Range("H138").Value = 2827572
$ws.Range("I138").Value = 1016.5263
$ws.Range("J138").Value = 4170185.8
$ws.Range("K138").Value = 3049.5789
$ws.Range("L138").Value = 12510557.4
$ws.Range("M138").Value = 2090.4211
$ws.Range("N138").Value = -12520837.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1711.38
$ws.Range("I32").Value = 1446.6967
$ws.Range("J32").Value = 3852.9092
$ws.Range("K32").Value = 1446.6967
$ws.Range("L32").Value = 3852.9092
$ws.Range("M32").Value = -1159.6967
$ws.Range("N32").Value = -4426.9092
$ws.Range("H61").Value = 111334830
$ws.Range("I61").Value = 166834080
$ws.Range("J61").Value = 336333.34
$ws.Range("K61").Value = 166834080
$ws.Range("L61").Value = 336333.34
$ws.Range("M61").Value = -166833868
$ws.Range("N61").Value = -336757.34
$ws.Range("H74").Value = 7002158
$ws.Range("I74").Value = 10914756
$ws.Range("J74").Value = 79869.234
$ws.Range("K74").Value = 10914756
$ws.Range("L74").Value = 79869.234
$ws.Range("M74").Value = -10913882
$ws.Range("N74").Value = -81617.234
$ws.Range("H77").Value = 7002158
$ws.Range("I77").Value = 10914756
$ws.Range("J77").Value = 79869.234
$ws.Range("K77").Value = 54573780
$ws.Range("L77").Value = 399346.17
$ws.Range("M77").Value = -54569412
$ws.Range("N77").Value = -408082.17
$ws.Range("H132").Value = 32816.254
$ws.Range("I132").Value = 22942.979
$ws.Range("K132").Value = 68828.93700000001
$ws.Range("M132").Value = -66298.93700000001
$ws.Range("H136").Value = 111334830
$ws.Range("I136").Value = 166834080
$ws.Range("J136").Value = 336333.34
$ws.Range("K136").Value = 500502240
$ws.Range("L136").Value = 1009000.02
$ws.Range("M136").Value = -500499690
$ws.Range("N136").Value = -1014100.02

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 585.61536
$ws.Range("I94").Value = 598.1111
$ws.Range("J94").Value = 557.5
$ws.Range("K94").Value = 598.1111
$ws.Range("L94").Value = 557.5
$ws.Range("M94").Value = -147.1111
$ws.Range("N94").Value = -1459.5
$ws.Range("H107").Value = 2119.111
$ws.Range("I107").Value = 2233.875
$ws.Range("J107").Value = 1201
$ws.Range("K107").Value = 2233.875
$ws.Range("L107").Value = 1201
$ws.Range("M107").Value = -313.875
$ws.Range("N107").Value = -5041
$ws.Range("H134").Value = 2533.56
$ws.Range("I134").Value = 1997.025
$ws.Range("J134").Value = 4679.7
$ws.Range("K134").Value = 5991.075000000001
$ws.Range("L134").Value = 14039.1
$ws.Range("M134").Value = -3456.075000000001
$ws.Range("N134").Value = -19109.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 927.8570999999999
$ws.Range("I105").Value = 1090
$ws.Range("J105").Value = 636
$ws.Range("K105").Value = 1090
$ws.Range("L105").Value = 636
$ws.Range("M105").Value = 657
$ws.Range("N105").Value = -4130
$ws.Range("H107").Value = 386.55173
$ws.Range("I107").Value = 371.57895
$ws.Range("J107").Value = 415
$ws.Range("K107").Value = 371.57895
$ws.Range("L107").Value = 415
$ws.Range("M107").Value = 1548.42105
$ws.Range("N107").Value = -4255
$ws.Range("H134").Value = 25913.299
$ws.Range("I134").Value = 1998.0834
$ws.Range("J134").Value = 104181.27
$ws.Range("K134").Value = 5994.2502
$ws.Range("L134").Value = 312543.81
$ws.Range("M134").Value = -3459.2502
$ws.Range("N134").Value = -317613.81

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1027.7261
$ws.Range("I131").Value = 736
$ws.Range("J131").Value = 1049.1765
$ws.Range("K131").Value = 2208
$ws.Range("L131").Value = 3147.5295
$ws.Range("M131").Value = 2832
$ws.Range("N131").Value = -13227.5295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37875.03
$ws.Range("I70").Value = 58027.684
$ws.Range("K70").Value = 58027.684
$ws.Range("M70").Value = -57757.684
$ws.Range("H73").Value = 37875.03
$ws.Range("I73").Value = 58027.684
$ws.Range("K73").Value = 58027.684
$ws.Range("M73").Value = -57091.684
$ws.Range("H132").Value = 39841.94
$ws.Range("I132").Value = 26796.129
$ws.Range("J132").Value = 78979.38
$ws.Range("K132").Value = 80388.387
$ws.Range("L132").Value = 236938.14
$ws.Range("M132").Value = -77858.387
$ws.Range("N132").Value = -241998.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1857.1428
$ws.Range("H71").Value = 1857.1428
$ws.Range("H136").Value = 69602.836
$ws.Range("I136").Value = 41419.4
$ws.Range("J136").Value = 210520
$ws.Range("K136").Value = 124258.2
$ws.Range("L136").Value = 631560
$ws.Range("M136").Value = -121708.2
$ws.Range("N136").Value = -636660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 30776.666
$ws.Range("J94").Value = 30776.666
$ws.Range("L94").Value = 30776.666
$ws.Range("N94").Value = -32578.666
$ws.Range("H107").Value = 385.53333
$ws.Range("I107").Value = 292.33334
$ws.Range("J107").Value = 447.66666
$ws.Range("K107").Value = 877.0000200000001
$ws.Range("L107").Value = 1342.99998
$ws.Range("M107").Value = 1042.99998
$ws.Range("N107").Value = -5182.999980000001
$ws.Range("H122").Value = 2669.1428
$ws.Range("I122").Value = 2019.4615
$ws.Range("K122").Value = 6058.3845
$ws.Range("M122").Value = -3608.3845
$ws.Range("H136").Value = 44003.086
$ws.Range("I136").Value = 24352.209
$ws.Range("J136").Value = 255250
$ws.Range("K136").Value = 73056.62699999999
$ws.Range("L136").Value = 765750
$ws.Range("M136").Value = -70506.62699999999
$ws.Range("N136").Value = -770850
